$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Pv primer sets" sheet: recolor the four "Name" cells red (A2, A5, A8,
#    A11) - order matters so the new font/fill xf combos get minted in the
#    same order as the target file (fill 2, 3, 4, 5, then 0 on the other
#    sheet).  Also narrow/selection changes.
# ---------------------------------------------------------------------------
$wsSets = $wb.Worksheets.Item("Pv primer sets")

$wsSets.Range("A8").Font.Color = 255
$wsSets.Range("A11").Font.Color = 255
$wsSets.Range("A5").Font.Color = 255
$wsSets.Range("A2").Font.Color = 255

$wsSets.Columns.Item(1).ColumnWidth = 18.167

# ---------------------------------------------------------------------------
# 2) "Pv primer order" sheet: fill in the primer-ordering table.
# ---------------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Item("Pv primer order")

$wsOrder.Columns.Item(1).ColumnWidth = 20.667

$names = @(
  "Pf_Hu_P1","Pf_Hu_P2","Pf_Hu_P3","Pf_Hu_P4","Pf_Hu_P5","Pf_Hu_P6","Pf_Hu_P7",
  "Pf_Hu_P8","Pf_Hu_P9","Pf_Hu_P10","Pf_Hu_P11","Pf_Hu_P12","Pf_Hu_P13","Pf_Hu_P14",
  "Pf_Hu_P15","Pf_Hu_P16","Pf_Hu_P17","Pf_Hu_P18"
)
$seqs = @(
  "CGAAATAT","CGAATAAA","CGAATAAT","CGATAAAA","CGTAAATA","TAATCGTA","ATTTTTTACG",
  "CGAAATGTA","CGAATACG","CGATAACG","CGATTACG","CTTTTACGA","GACGAAATA","TACGAATTG",
  "TTATGTACG","TTTTTTTACG","TATAACGA","TATTTTTACG"
)

$wsOrder.Cells.Item(1,1).Value = "Primer Name"
$wsOrder.Cells.Item(1,2).Value = "Primers to Order"
$wsOrder.Cells.Item(1,1).Font.Bold = $true
$wsOrder.Cells.Item(1,2).Font.Bold = $true

for ($i = 0; $i -lt $names.Length; $i++) {
  $r = $i + 2
  $wsOrder.Cells.Item($r,1).Value = $names[$i]
  $wsOrder.Cells.Item($r,2).Value = $seqs[$i]
}

# Red (no fill) highlight for the unique primers: rows 2-7 and row 18.
$redRows = @(2,3,4,5,6,7,18)
foreach ($r in $redRows) {
  $wsOrder.Cells.Item($r,1).Font.Color = 255
  $wsOrder.Cells.Item($r,2).Font.Color = 255
}

# Notes in column D.
$wsOrder.Cells.Item(1,4).Value = "Order primer sets 10 and 13 which had no overlaps"
$wsOrder.Cells.Item(2,4).Value = "Order primer sets 1 and 4 had no overlaps"
$wsOrder.Cells.Item(4,4).Value = "All together, there were 9 duplicates"
$wsOrder.Cells.Item(7,4).Value = "Red primers are unique from order Pv_HuPf Primer sets"
$wsOrder.Cells.Item(7,4).Font.Color = 255

# Make sure "Pv primer sets" shows its own selection/active cell first...
$wsSets.Range("A2").Select()

# ...then leave "Pv primer order" as the active/selected sheet+range.
$wsOrder.Activate()
$wsOrder.Range("A2:A7").Select()
